# Changes made for Telemundo runs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Swap the TcNo values in A2 and A3
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

# Update the active selection on the TestData sheet
$ws.Activate()
$ws.Range("A3").Select()
